$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("EMP001",    "Maruthi.M",   "02:07 am", "01:00 pm", "10.87"),
    @("EMP575206", "Bhavya",      "08:43 am", "01:00 pm", "4.27"),
    @("EMP656501", "Vrushvini",   "09:58 am", "10:21 am", "0.39"),
    @("EMP940311", "Om Prakash",  "09:59 am", "01:00 pm", "3.00"),
    @("EMP084414", "Chinmayee",   "10:12 am", "10:13 am", "0.02"),
    @("EMP744311", "Sneha",       "10:05 am", "01:00 pm", "2.90"),
    @("EMP025103", "Karthik",     "10:09 am", "01:00 pm", "2.84"),
    @("EMP995605", "Akhil",       "10:10 am", "01:00 pm", "2.83")
)

# Column E ("Working Hours") holds numeric-looking strings (e.g. "10.87").
# Force that range to be treated as plain text so Excel does not
# auto-convert the values into real numbers.
$fillRange = $ws.Range("E2:E9")
$fillRange.NumberFormat = "@"

$row = 2
foreach ($entry in $data) {
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $ws.Cells.Item($row, 3).Value = $entry[2]
    $ws.Cells.Item($row, 4).Value = $entry[3]
    $ws.Cells.Item($row, 5).Value = $entry[4]
    $row++
}

# Restore default (unstyled) cell appearance now that the text values
# have been committed as strings.
$fillRange.Style = "Normal"

$ws.Columns.Item(3).ColumnWidth = 17.1
$ws.Columns.Item(4).ColumnWidth = 17.1
$ws.Columns.Item(5).ColumnWidth = 15.1
